$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEGFA165_NRP1")
$ws.Range("A20").Value = "a"
$ws.Range("B20").Value = "b"
$ws.Range("C20").Value = "c"
$ws.Range("D20").Value = "d"
$ws.Range("E20").Value = "e"

$ws.Range("A20").Interior.PatternColorIndex = -4105
$ws.Range("B20").Interior.TintAndShade = 0
$ws.Range("C20").Interior.PatternTintAndShade = 0
$ws.Range("D20").Interior.PatternColor = 0
$ws.Range("E20").Interior.Pattern = 0
